$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(13.45402668174342,42.825613359375,80.71564164884869,46.5472306737842,65.4254602287329,84.30368978368162,103.18191935625,122.0601489288184,140.9383784837671,151.6360419)
for ($r = 6; $r -le 10; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(12.16610275813231,38.726013046875,72.9889135585709,42.09136824830228,59.16242705397386,76.23348585964544,93.30454468124998,110.3756035028546,127.4466623085261,137.1202623)
for ($r = 11; $r -le 11; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(21.83164398985746,69.492469921875,130.9760411938048,75.53148160191893,106.1648969021655,136.7983122024122,167.43172753125,198.0651428600879,228.6985581603344,246.0574935)
for ($r = 12; $r -le 12; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(13.97630816931652,44.48809148437499,83.84899989080775,48.35418092397961,67.96525797863518,87.57633503329075,107.18741210625,126.7984891792092,146.4095662338647,157.5225099)
for ($r = 37; $r -le 54; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(16.49159358872442,52.494515390625,98.93911984968933,57.05637643742079,80.1968158656646,103.3372552939084,126.47769474375,149.6181341935916,172.7585736218354,185.8714893)
for ($r = 57; $r -le 61; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(14.23467082008407,45.310487578125,85.39901221262792,49.2480446115763,69.22164729108685,89.19524997059742,109.16885266875,129.1424553669026,149.1160580464131,160.4344329)
for ($r = 62; $r -le 83; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(19.10966844983553,60.828129140625,114.6459113774671,66.11419513840042,92.9282275651748,119.7422599919492,146.55629244375,173.3703248955509,200.1843573223251,215.3789757)
for ($r = 88; $r -le 89; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(18.30874423245614,58.27870125,109.8408731798245,63.34321770685072,89.03342069657464,114.7236236862986,140.4138267,166.1040297137015,191.7942327034253,206.3520144)
for ($r = 105; $r -le 143; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(16.65716793265716,53.02155632812499,99.93246107958699,57.62921811248273,81.00198579063363,104.3747534687845,127.74752116875,151.1202888687155,174.4930565468663,187.7376249)
for ($r = 144; $r -le 148; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}

$vals = @(23.68074270330775,75.37834992187501,142.0694627142726,81.92885440261081,115.1568617018196,148.3848690010284,181.61287633125,214.8408836614716,248.0688909606803,266.8980951)
for ($r = 149; $r -le 149; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c-1]
    }
}
